$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.131.41"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "3.945.39"
$ws.Range("E3").Value = "  -3.14%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.11"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.14"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "3.937.12"
$ws.Range("E7").Value = "  -3.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.685"
$ws.Range("E8").Value = "  -5.97%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("E11").Value = "  -6.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.76"
$ws.Range("E12").Value = "  +12.18%  "
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.66"
$ws.Range("E14").Value = "  -4.08%  "
$ws.Range("D15").Value = "4.578.96"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "3.950.20"
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.97"
$ws.Range("E17").Value = "  -4.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.50"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("E20").Value = "  -6.23%  "
$ws.Range("D21").Value = "71.010.23"
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "427.18"
$ws.Range("E22").Value = "  -4.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "97.26"
$ws.Range("E23").Value = "  -6.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.57"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.23"
$ws.Range("E25").Value = "  +4.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.58"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.32"
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.92"
$ws.Range("E28").Value = "  +19.15%  "
$ws.Range("E29").Value = "  -3.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.91"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.54"
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.85"
$ws.Range("E32").Value = "  +18.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.71"
$ws.Range("E33").Value = "  +19.64%  "
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.35"
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "681.16"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("D39").Value = "0.0₃0814"
$ws.Range("E39").Value = "  -6.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.149"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.37"
$ws.Range("E41").Value = "  -3.06%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0483"
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.26"
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("E47").Value = "  -6.42%  "
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.35"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.65"
$ws.Range("E51").Value = "  +0.57%  "
